$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-01-11 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-12 Monday", 2) | Out-Null
$d.Content.Find.Execute("61÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "82÷7=", 2) | Out-Null
$d.Content.Find.Execute("39÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "51÷2=", 2) | Out-Null
$d.Content.Find.Execute("32÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "16÷4=", 2) | Out-Null
$d.Content.Find.Execute("55÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "49÷9=", 2) | Out-Null
$d.Content.Find.Execute("73÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "86÷8=", 2) | Out-Null
$d.Content.Find.Execute("29÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷7=", 2) | Out-Null
$d.Content.Find.Execute("40÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "25÷4=", 2) | Out-Null
$d.Content.Find.Execute("96÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "63÷7=", 2) | Out-Null
$d.Content.Find.Execute("54÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "74÷6=", 2) | Out-Null
$d.Content.Find.Execute("37÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "61÷2=", 2) | Out-Null
$d.Content.Find.Execute("93÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "61÷4=", 2) | Out-Null
$d.Content.Find.Execute("76÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "98÷8=", 2) | Out-Null
$d.Content.Find.Execute("28÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "56÷5=", 2) | Out-Null
$d.Content.Find.Execute("94÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "48÷8=", 2) | Out-Null
$d.Content.Find.Execute("86÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "33÷7=", 2) | Out-Null
$d.Content.Find.Execute("15÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "48÷5=", 2) | Out-Null
$d.Content.Find.Execute("38÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "66÷2=", 2) | Out-Null
$d.Content.Find.Execute("43÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "86÷6=", 2) | Out-Null
$d.Content.Find.Execute("61÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "31÷3=", 2) | Out-Null
$d.Content.Find.Execute("29÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "92÷5=", 2) | Out-Null
$d.Content.Find.Execute("32÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "67÷8=", 2) | Out-Null
$d.Content.Find.Execute("67÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "38÷7=", 2) | Out-Null
$d.Content.Find.Execute("26÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "79÷8=", 2) | Out-Null
$d.Content.Find.Execute("59÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "74÷4=", 2) | Out-Null
$d.Content.Find.Execute("78÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "99÷8=", 2) | Out-Null
